$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G (7th column), shifting "Ciudad de origen",
# "Departamento/Estado de origen" and "Pais de origen" one column to the right.
$ws.Columns.Item(7).Insert()

# New header + value for the inserted "Genero" column.
$ws.Cells.Item(1, 7).Value = "Género"
$ws.Cells.Item(2, 7).Value = "masculino"

# Column B width: 7.75 -> 8.0
$ws.Columns.Item(2).ColumnWidth = 7.15

# New column G should match column F's width (10.88)
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(6).ColumnWidth
